# Update calculated result values for pl_mw.xlsx (Case_3_234, line results)
# Commit: "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.044143055898587
$ws.Range("C2").Value = 0.06022667446837104
$ws.Range("D2").Value = 0.3270776345368063
$ws.Range("F2").Value = 4.778112276208446
$ws.Range("G2").Value = 0.002614359305205317
$ws.Range("K2").Value = 0.5438196910038471
$ws.Range("L2").Value = 0.2628897277578517
$ws.Range("M2").Value = 0.2749282348159738

# Row 3
$ws.Range("B3").Value = 1.028915416720679
$ws.Range("C3").Value = 0.05801084374200727
$ws.Range("D3").Value = 0.3144497402019226
$ws.Range("F3").Value = 4.564674969141379
$ws.Range("G3").Value = 0.002619346392822724
$ws.Range("K3").Value = 0.530062495270613
$ws.Range("L3").Value = 0.2589639369865324
$ws.Range("M3").Value = 0.2709598082021856

# Row 4
$ws.Range("B4").Value = 1.020502663823834
$ws.Range("C4").Value = 0.05660976884723112
$ws.Range("D4").Value = 0.3066083428495858
$ws.Range("F4").Value = 4.43378742677848
$ws.Range("G4").Value = 0.002622566071081746
$ws.Range("K4").Value = 0.5222778340966983
$ws.Range("L4").Value = 0.2566710397393948
$ws.Range("M4").Value = 0.2687222581441588

# Row 5
$ws.Range("B5").Value = 1.017309622114709
$ws.Range("C5").Value = 0.05602844316561928
$ws.Range("D5").Value = 0.3033904164351924
$ws.Range("F5").Value = 4.380486876569421
$ws.Range("G5").Value = 0.00262391788221197
$ws.Range("K5").Value = 0.5192714361903086
$ws.Range("L5").Value = 0.2557661997502336
$ws.Range("M5").Value = 0.2678604371157043

# Row 6
$ws.Range("B6").Value = 1.016793615113357
$ws.Range("C6").Value = 0.05593128250278312
$ws.Range("D6").Value = 0.3028547116848728
$ws.Range("F6").Value = 4.371638490438698
$ws.Range("G6").Value = 0.002624144755269342
$ws.Range("K6").Value = 0.5187822287900445
$ws.Range("L6").Value = 0.2556177348442077
$ws.Range("M6").Value = 0.2677203504595163

# Row 7
$ws.Range("B7").Value = 1.020458649480332
$ws.Range("C7").Value = 0.05660197113562049
$ws.Range("D7").Value = 0.3065650363785295
$ws.Range("F7").Value = 4.433068451652304
$ws.Range("G7").Value = 0.002622584140959106
$ws.Range("K7").Value = 0.5222366178103215
$ws.Range("L7").Value = 0.2566587171994854
$ws.Range("M7").Value = 0.2687104329613597

# Row 8
$ws.Range("B8").Value = 1.03869781357858
$ws.Range("C8").Value = 0.05947096216314307
$ws.Range("D8").Value = 0.322741513323578
$ws.Range("F8").Value = 4.7044826052269
$ws.Range("G8").Value = 0.002616046230353132
$ws.Range("K8").Value = 0.5389383579275062
$ws.Range("L8").Value = 0.2615116990359354
$ws.Range("M8").Value = 0.2735185622772498

# Row 9
$ws.Range("B9").Value = 1.081923360638257
$ws.Range("C9").Value = 0.06478356804336372
$ws.Range("D9").Value = 0.3537863677284037
$ws.Range("F9").Value = 5.238232127359481
$ws.Range("G9").Value = 0.00260446939393711
$ws.Range("K9").Value = 0.5769775059563642
$ws.Range("L9").Value = 0.2719631277731622
$ws.Range("M9").Value = 0.2845307687727967

# Row 10
$ws.Range("B10").Value = 1.118267180026351
$ws.Range("C10").Value = 0.06850716210200147
$ws.Range("D10").Value = 0.3762110059565202
$ws.Range("F10").Value = 5.631621853487246
$ws.Range("G10").Value = 0.002596713283677366
$ws.Range("K10").Value = 0.6081973809864678
$ws.Range("L10").Value = 0.2802156436616627
$ws.Range("M10").Value = 0.2935937150799006

# Row 11
$ws.Range("B11").Value = 1.135805831765424
$ws.Range("C11").Value = 0.07016462934291923
$ws.Range("D11").Value = 0.3863354803887376
$ws.Range("F11").Value = 5.810928109450572
$ws.Range("G11").Value = 0.002593345632172284
$ws.Range("K11").Value = 0.6231218838657071
$ws.Range("L11").Value = 0.2840954999655168
$ws.Range("M11").Value = 0.2979294529345395

# Row 12
$ws.Range("B12").Value = 1.14259252223431
$ws.Range("C12").Value = 0.07078723403239096
$ws.Range("D12").Value = 0.3901588148048347
$ws.Range("F12").Value = 5.878882407021194
$ws.Range("G12").Value = 0.002592093343602433
$ws.Range("K12").Value = 0.6288781283609239
$ws.Range("L12").Value = 0.2855828401754223
$ws.Range("M12").Value = 0.2996020173551273

# Row 13
$ws.Range("B13").Value = 1.14112442065607
$ws.Range("C13").Value = 0.07065336611861284
$ws.Range("D13").Value = 0.3893358551744939
$ws.Range("F13").Value = 5.864244720440126
$ws.Range("G13").Value = 0.002592362026974328
$ws.Range("K13").Value = 0.6276337508508334
$ws.Range("L13").Value = 0.2852617081619542
$ws.Range("M13").Value = 0.2992404331469629

# Row 14
$ws.Range("B14").Value = 1.136361263847562
$ws.Range("C14").Value = 0.07021595134076719
$ws.Range("D14").Value = 0.3866502380567454
$ws.Range("F14").Value = 5.81651762455391
$ws.Range("G14").Value = 0.002593242146158802
$ws.Range("K14").Value = 0.6235933515402223
$ws.Range("L14").Value = 0.2842175008911454
$ws.Range("M14").Value = 0.2980664397133026

# Row 15
$ws.Range("B15").Value = 1.133462617199484
$ws.Range("C15").Value = 0.06994737130633766
$ws.Range("D15").Value = 0.3850038544867971
$ws.Range("F15").Value = 5.787290682540856
$ws.Range("G15").Value = 0.002593784231768284
$ws.Range("K15").Value = 0.6211321432989223
$ws.Range("L15").Value = 0.283580255157716
$ws.Range("M15").Value = 0.297351337558581

# Row 16
$ws.Range("B16").Value = 1.117141272788473
$ws.Range("C16").Value = 0.06839812653347366
$ws.Range("D16").Value = 0.375547843066073
$ws.Range("F16").Value = 5.619911211206897
$ws.Range("G16").Value = 0.00259693658881998
$ws.Range("K16").Value = 0.6072366325692826
$ws.Range("L16").Value = 0.2799646196908867
$ws.Range("M16").Value = 0.2933146574155927

# Row 17
$ws.Range("B17").Value = 1.107386631837329
$ws.Range("C17").Value = 0.0674385261452386
$ws.Range("D17").Value = 0.3697275719801212
$ws.Range("F17").Value = 5.517322410133318
$ws.Range("G17").Value = 0.002598911508012567
$ws.Range("K17").Value = 0.598897763726967
$ws.Range("L17").Value = 0.2777787792854127
$ws.Range("M17").Value = 0.2908928925530176

# Row 18
$ws.Range("B18").Value = 1.101870643240176
$ws.Range("C18").Value = 0.06688315019242097
$ws.Range("D18").Value = 0.3663726695237557
$ws.Range("F18").Value = 5.458348844485386
$ws.Range("G18").Value = 0.002600062557407275
$ws.Range("K18").Value = 0.5941694176621013
$ws.Range("L18").Value = 0.2765333748987047
$ws.Range("M18").Value = 0.2895199964024044

# Row 19
$ws.Range("B19").Value = 1.10001925875946
$ws.Range("C19").Value = 0.0666945128623837
$ws.Range("D19").Value = 0.3652355012227275
$ws.Range("F19").Value = 5.438386908868637
$ws.Range("G19").Value = 0.00260045488532778
$ws.Range("K19").Value = 0.5925801281685779
$ws.Range("L19").Value = 0.2761137335239567
$ws.Range("M19").Value = 0.2890585963321186

# Row 20
$ws.Range("B20").Value = 1.108415232857936
$ws.Range("C20").Value = 0.06754103178256798
$ws.Range("D20").Value = 0.370347895588452
$ws.Range("F20").Value = 5.528239736895699
$ws.Range("G20").Value = 0.002598699709681886
$ws.Range("K20").Value = 0.5997784123568408
$ws.Range("L20").Value = 0.2780102407456582
$ws.Range("M20").Value = 0.291148618953585

# Row 21
$ws.Range("B21").Value = 1.137756373585603
$ws.Range("C21").Value = 0.07034456595298977
$ws.Range("D21").Value = 0.387439352709066
$ws.Range("F21").Value = 5.830534703984938
$ws.Range("G21").Value = 0.002592983011473009
$ws.Range("K21").Value = 0.624777268078816
$ws.Range("L21").Value = 0.2845237176952651
$ws.Range("M21").Value = 0.2984104357110482

# Row 22
$ws.Range("B22").Value = 1.157779024354824
$ws.Range("C22").Value = 0.07214753906175631
$ws.Range("D22").Value = 0.398548116515542
$ws.Range("F22").Value = 6.028424078739818
$ws.Range("G22").Value = 0.002589380632739614
$ws.Range("K22").Value = 0.6417258706167672
$ws.Range("L22").Value = 0.2888863054438247
$ws.Range("M22").Value = 0.3033355258878245

# Row 23
$ws.Range("B23").Value = 1.147014909165961
$ws.Range("C23").Value = 0.07118787466393428
$ws.Range("D23").Value = 0.3926246434395466
$ws.Range("F23").Value = 5.92277585211491
$ws.Range("G23").Value = 0.002591291089561605
$ws.Range("K23").Value = 0.6326239699095026
$ws.Range("L23").Value = 0.2865482289287939
$ws.Range("M23").Value = 0.3006904962181238

# Row 24
$ws.Range("B24").Value = 1.107949915547039
$ws.Range("C24").Value = 0.06749470049073381
$ws.Range("D24").Value = 0.3700674745829815
$ws.Range("F24").Value = 5.52330399478538
$ws.Range("G24").Value = 0.002598795414786118
$ws.Range("K24").Value = 0.5993800662707827
$ws.Range("L24").Value = 0.2779055619372741
$ws.Range("M24").Value = 0.2910329445995643

# Row 25
$ws.Range("B25").Value = 1.069427236145884
$ws.Range("C25").Value = 0.06337874044557523
$ws.Range("D25").Value = 0.3454574379556021
$ws.Range("F25").Value = 5.093643978159804
$ws.Range("G25").Value = 0.002607468987460441
$ws.Range("K25").Value = 0.5661158426975703
$ws.Range("L25").Value = 0.2690352610656959
$ws.Range("M25").Value = 0.2813814469509879
